# experiment template new experiments
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "experiment_description" (sheet1) - add experiments 7 and 8 just
# before the existing 1001 block.
# ---------------------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("experiment_description")
$wsDesc.Range("A8:A9").EntireRow.Insert()

$wsDesc.Range("A8").Value = 7
$wsDesc.Range("B8").Value = "LOS driven:3 + treatment constraints splitting of transitions  in Inpatient Ward"
$wsDesc.Range("C8").Value = "base"
$wsDesc.Range("D8").Value = "1;4"

$wsDesc.Range("A9").Value = 8
$wsDesc.Range("B9").Value = "LOS driven:7 + treatment_constraints splitting in length of stay in Inpatient Ward"
$wsDesc.Range("C9").Value = "base"
$wsDesc.Range("D9").Value = "1;4"

# ---------------------------------------------------------------------------
# Sheet "experiment_specification" (sheet2) - add the state rows describing
# experiments 7 and 8, just before the existing 1001 block.
# ---------------------------------------------------------------------------
$wsSpec = $wb.Worksheets.Item("experiment_specification")
$wsSpec.Range("A20:A25").EntireRow.Insert()

$wsSpec.Range("A20").Value = 7
$wsSpec.Range("B20").Value = "home"
$wsSpec.Range("C20").Value = "length_of_stay_simple_two_weeks"
$wsSpec.Range("D20").Value = "age_simple"
$wsSpec.Range("E20").Value = "age_simple"

$wsSpec.Range("A21").Value = 7
$wsSpec.Range("B21").Value = "inpatient_ward"
$wsSpec.Range("C21").Value = "none"
$wsSpec.Range("D21").Value = "age_simple_intensive_care_unit_restriction"
$wsSpec.Range("E21").Value = "none"

$wsSpec.Range("A22").Value = 7
$wsSpec.Range("B22").Value = "intensive_care_unit"
$wsSpec.Range("C22").Value = "none"
$wsSpec.Range("D22").Value = "age_simple"
$wsSpec.Range("E22").Value = "none"

$wsSpec.Range("A23").Value = 8
$wsSpec.Range("B23").Value = "home"
$wsSpec.Range("C23").Value = "length_of_stay_simple_two_weeks"
$wsSpec.Range("D23").Value = "age_simple"
$wsSpec.Range("E23").Value = "age_simple"

$wsSpec.Range("A24").Value = 8
$wsSpec.Range("B24").Value = "inpatient_ward"
$wsSpec.Range("C24").Value = "none"
$wsSpec.Range("D24").Value = "age_simple_intensive_care_unit_restriction"
$wsSpec.Range("E24").Value = "age_simple_intensive_care_unit_restriction"

$wsSpec.Range("A25").Value = 8
$wsSpec.Range("B25").Value = "intensive_care_unit"
$wsSpec.Range("C25").Value = "none"
$wsSpec.Range("D25").Value = "age_simple"
$wsSpec.Range("E25").Value = "none"

# ---------------------------------------------------------------------------
# Sheet "run_description" (sheet3) - append new run 7.
# ---------------------------------------------------------------------------
$wsRunDesc = $wb.Worksheets.Item("run_description")
$wsRunDesc.Range("A8").Value = 7
$wsRunDesc.Range("B8").Value = "Testing treatment constraint splitting"

# ---------------------------------------------------------------------------
# Sheet "run_specification" (sheet4) - append run 7 -> experiments 3, 7, 8.
# ---------------------------------------------------------------------------
$wsRunSpec = $wb.Worksheets.Item("run_specification")
$wsRunSpec.Range("A18").Value = 7
$wsRunSpec.Range("B18").Value = 3

$wsRunSpec.Range("A19").Value = 7
$wsRunSpec.Range("B19").Value = 7

$wsRunSpec.Range("A20").Value = 7
$wsRunSpec.Range("B20").Value = 8

# ---------------------------------------------------------------------------
# Restore per-sheet selections (last selection remembered per worksheet);
# the final active tab follows the last .Select() call below, which keeps
# "run_specification" active (same as before editing).
# ---------------------------------------------------------------------------
$wsDesc.Range("D9").Select()
$wsSpec.Range("D30").Select()
$wsRunDesc.Range("B18").Select()
$wsRunSpec.Range("B20").Select()
